$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Step 1: direct value edits, using the ORIGINAL row numbering
#     (table structure is still untouched at this point) ---
$t.Cell(1, 1).Range.Text  = "0M"      # was 99.97
$t.Cell(2, 1).Range.Text  = "0M"      # was 0.06
$t.Cell(3, 1).Range.Text  = "0M"      # was 227
$t.Cell(4, 1).Range.Text  = "340"     # was 85
$t.Cell(6, 1).Range.Text  = "0.00066" # was 0.00013
$t.Cell(8, 1).Range.Text  = "0.00019" # was 0.00002
$t.Cell(10, 1).Range.Text = "0.00029" # was 0.00006
$t.Cell(11, 1).Range.Text = "0.00032" # was 0.00009
$t.Cell(12, 1).Range.Text = "0.00034" # was 0.00589

# --- Step 2: delete the row that held 0.00007 (original row 7) ---
$t.Rows(7).Delete()

# --- Step 3: insert a new row just before the row that now holds "100.0"
#     (originally row 13, now row 12 after the deletion above) and set
#     its value ---
$newRow = $t.Rows.Add($t.Rows(12))
$newRow.Cells(1).Range.Text = "0.06399"

# --- Step 4: collapse the tab-separated runs in the final 3 rows down to
#     a single value (net row count unaffected: -1 deletion + 1 insertion) ---
$t.Cell(44, 1).Range.Text = "99.97"
$t.Cell(45, 1).Range.Text = "0.06"
$t.Cell(46, 1).Range.Text = "227"
